$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-24 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("126÷7=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "956÷3=318, 2", 2) | Out-Null
$d.Content.Find.Execute("540÷7=77, 1", $true, $false, $false, $false, $false, $true, 1, $false, "911÷9=101, 2", 2) | Out-Null
$d.Content.Find.Execute("737÷6=122, 5", $true, $false, $false, $false, $false, $true, 1, $false, "470÷8=58, 6", 2) | Out-Null
$d.Content.Find.Execute("450÷9=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "211÷8=26, 3", 2) | Out-Null
$d.Content.Find.Execute("788÷4=197, 0", $true, $false, $false, $false, $false, $true, 1, $false, "446÷6=74, 2", 2) | Out-Null
$d.Content.Find.Execute("751÷3=250, 1", $true, $false, $false, $false, $false, $true, 1, $false, "874÷5=174, 4", 2) | Out-Null
$d.Content.Find.Execute("977÷9=108, 5", $true, $false, $false, $false, $false, $true, 1, $false, "960÷8=120, 0", 2) | Out-Null
$d.Content.Find.Execute("244÷5=48, 4", $true, $false, $false, $false, $false, $true, 1, $false, "359÷4=89, 3", 2) | Out-Null
$d.Content.Find.Execute("897÷7=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "676÷5=135, 1", 2) | Out-Null
$d.Content.Find.Execute("552÷8=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "151÷9=16, 7", 2) | Out-Null
$d.Content.Find.Execute("904÷5=180, 4", $true, $false, $false, $false, $false, $true, 1, $false, "268÷5=53, 3", 2) | Out-Null
$d.Content.Find.Execute("755÷2=377, 1", $true, $false, $false, $false, $false, $true, 1, $false, "913÷9=101, 4", 2) | Out-Null
$d.Content.Find.Execute("608÷9=67, 5", $true, $false, $false, $false, $false, $true, 1, $false, "959÷8=119, 7", 2) | Out-Null
$d.Content.Find.Execute("728÷7=104, 0", $true, $false, $false, $false, $false, $true, 1, $false, "964÷5=192, 4", 2) | Out-Null
$d.Content.Find.Execute("114÷6=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "288÷9=32, 0", 2) | Out-Null
$d.Content.Find.Execute("859÷5=171, 4", $true, $false, $false, $false, $false, $true, 1, $false, "445÷7=63, 4", 2) | Out-Null
$d.Content.Find.Execute("251÷4=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "929÷9=103, 2", 2) | Out-Null
$d.Content.Find.Execute("858÷9=95, 3", $true, $false, $false, $false, $false, $true, 1, $false, "772÷6=128, 4", 2) | Out-Null
$d.Content.Find.Execute("627÷5=125, 2", $true, $false, $false, $false, $false, $true, 1, $false, "234÷4=58, 2", 2) | Out-Null
$d.Content.Find.Execute("846÷2=423, 0", $true, $false, $false, $false, $false, $true, 1, $false, "273÷4=68, 1", 2) | Out-Null
$d.Content.Find.Execute("958÷4=239, 2", $true, $false, $false, $false, $false, $true, 1, $false, "649÷5=129, 4", 2) | Out-Null
$d.Content.Find.Execute("661÷7=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "302÷3=100, 2", 2) | Out-Null
$d.Content.Find.Execute("568÷9=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "225÷7=32, 1", 2) | Out-Null
$d.Content.Find.Execute("966÷4=241, 2", $true, $false, $false, $false, $false, $true, 1, $false, "352÷8=44, 0", 2) | Out-Null
$d.Content.Find.Execute("107÷3=35, 2", $true, $false, $false, $false, $false, $true, 1, $false, "256÷6=42, 4", 2) | Out-Null
